$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking snapshot refresh: update prices / 1h-volume deltas, and the two
# swapped rank rows (ImmutableX <-> NEARProtocol).
#
# Cells whose new text looks like a plain number (e.g. "604.52") get a leading
# apostrophe so Excel stores them as text (matching the source inlineStr cells)
# instead of silently converting them to numeric values; the apostrophe itself
# is not stored in the cell, it only forces the text interpretation. The quote-
# prefix flag Excel attaches is then cleared via Style so no stray cell style
# is introduced.
$ws.Range("D2").Value = "66.951.46"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "3.517.17"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'604.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").Value = "'147.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.85%  "
$ws.Range("D7").Value = "3.515.49"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("D11").Value = "'7.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.06%  "
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("D14").Value = "4.114.26"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("D16").Value = "3.519.45"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "66.981.81"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("E19").Value = "  +7.99%  "
$ws.Range("D20").Value = "'6.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("E21").Value = "  -1.65%  "
$ws.Range("D22").Value = "'436.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.19%  "
$ws.Range("D23").Value = "'0.608"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.33%  "
$ws.Range("D24").Value = "'79.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("D25").Value = "3.658.09"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -4.05%  "
$ws.Range("D28").Value = "'9.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.91%  "
$ws.Range("D29").Value = "'8.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.13%  "
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("D35").Value = "3.514.83"
$ws.Range("E35").Value = "  +0.57%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'5.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.65%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.45%  "
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "'0.0890"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").Value = "'170.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.41%  "
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("E44").Value = "  -10.37%  "
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("E47").Value = "  +2.82%  "
$ws.Range("D48").Value = "'28.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.57%  "
$ws.Range("E49").Value = "  -1.90%  "
$ws.Range("D50").Value = "'2.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.74%  "
$ws.Range("D51").Value = "'0.989"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.10%  "
